# Adds a new "samples" worksheet (a small sample/material tracking table),
# fills in the "sample" column on Sheet1 (F2/F3) with the new sample IDs,
# and makes "samples" the active/selected sheet - matching the commit
# "added sample df merge to file list".

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")

# --- Fill in the sample IDs on Sheet1 (F2/F3 were blank, now "s1"/"s2") ---
$sheet1.Range("F2").Value = "s1"
$sheet1.Range("F3").Value = "s2"

# --- Add the new "samples" worksheet right after Sheet1 ---
$samples = $wb.Worksheets.Add($null, $sheet1)
$samples.Name = "samples"

# Header row
$samples.Range("A1").Value = "ID"
$samples.Range("B1").Value = "Owner "
$samples.Range("C1").Value = "Material"
$samples.Range("D1").Value = "Density"
$samples.Range("E1").Value = "Batch-date"
$samples.Range("F1").Value = "Solvent"
$samples.Range("G1").Value = "Concentration"
$samples.Range("H1").Value = "Dilution-date"
$samples.Range("I1").Value = "Prep-date"

# Style the header row: bold dark font + light fill, matching col_test1's other header row
$samples.Range("A1:I1").Font.Bold = $true
$samples.Range("A1:I1").Interior.ThemeColor = 8
$samples.Range("A1:I1").Interior.TintAndShade = 0.8
$samples.Range("A1:I1").WrapText = $true
$samples.Range("A1").HorizontalAlignment = -4152
$samples.Range("A1").NumberFormat = "@"
$samples.Range("B1:I1").NumberFormat = "@"
$samples.Range("G1:H1").NumberFormat = "0.00"
$samples.Rows.Item(1).RowHeight = 40

# Data rows
$samples.Range("A2").Value = "s1"
$samples.Range("B2").Value = "Lexie"
$samples.Range("C2").Value = "Iron oxide"
$samples.Range("F2").Value = "Water"
$samples.Range("G2").Value = 0.1

$samples.Range("A3").Value = "s2"
$samples.Range("B3").Value = "Arlo"
$samples.Range("C3").Value = "Dynabeads"
$samples.Range("F3").Value = "PBS"
$samples.Range("G3").Value = 0.2

# --- Make "samples" the active sheet / tab ---
$samples.Activate()
$samples.Range("C2").Select()
